# Applies the "Black River Gold" edit:
#   1. Removes the "Meta description: ..." paragraph that used to sit right
#      under the H1 title.
#   2. Inserts a new bold paragraph ("Play Black River Gold for Free -
#      Exciting Features & Max Win 5,000x") right before the closing
#      "Prompt: ..." paragraph at the end of the document.
#   3. Rewrites that closing paragraph's text to the old meta-description
#      copy (keeping its existing italic formatting).

function Get-ParagraphIndexAt($doc, $pos) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Start -le $pos -and $pos -lt $p.Range.End) {
            return $i
        }
    }
    return -1
}

function Find-ParagraphIndexByText($doc, $searchText) {
    $rng = $doc.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if (-not $found) { return -1 }
    return Get-ParagraphIndexAt $doc $rng.Start
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: delete the "Meta description" paragraph near the top.
# ---------------------------------------------------------------------
$metaIdx = Find-ParagraphIndexByText $d `
    "Meta description: Read our Black River Gold slot review and play for free. This ElK Studios game offers unique features, stacked symbols, and endless free spins."
if ($metaIdx -ne -1) {
    $d.Paragraphs.Item($metaIdx).Range.Delete()
}

# ---------------------------------------------------------------------
# Step 2: insert a new bold paragraph just above the trailing
# "Prompt: ..." paragraph, containing the old heading text.
# ---------------------------------------------------------------------
$headingText = "Play Black River Gold for Free - Exciting Features & Max Win 5,000x"

# Build the new paragraph by splitting a plain (Normal-style, non-bold,
# non-italic) paragraph so the inserted run doesn't inherit stray
# character formatting; then relocate it with Cut/Paste.
$anchorIdx = Find-ParagraphIndexByText $d "Do you hear that whistle blowin"
$anchorPara = $d.Paragraphs.Item($anchorIdx)
$splitPos = $anchorPara.Range.End - 1
$insertPoint = $d.Range($splitPos, $splitPos)
$insertPoint.InsertAfter("`r" + $headingText)

$newIdx = $anchorIdx + 1
$newPara = $d.Paragraphs.Item($newIdx)
$textOnly = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$textOnly.Font.Bold = $true

$newPara.Range.Cut() | Out-Null

$promptIdx = Find-ParagraphIndexByText $d "Prompt: Create a feature image"
$promptPara = $d.Paragraphs.Item($promptIdx)
$targetPos = $promptPara.Range.Start
$targetRng = $d.Range($targetPos, $targetPos)
$targetRng.Paste() | Out-Null

# ---------------------------------------------------------------------
# Step 3: rewrite the closing "Prompt: ..." paragraph's text, keeping
# its italic run formatting.
# ---------------------------------------------------------------------
$oldPrompt = "Prompt: Create a feature image for Black River Gold that captures the adventurous spirit of the game. The image should be in a cartoon style and feature a happy and confident Maya warrior wearing glasses. The warrior can be holding a revolver and standing in front of a mountain of gold and precious stones. The background should be set in the wild west with cacti and a clear blue sky. The overall image should be eye-catching and vibrant, with bold colors and dynamic lines to draw attention to the game's exciting features."
$newPrompt = "Read our Black River Gold slot review and play for free. This ElK Studios game offers unique features, stacked symbols, and endless free spins."

$rng = $d.Content
$rng.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newPrompt, 2) | Out-Null

Write-Output "done"
